# Registro_Actividades_Sistema - actualizacion de tiempos, Jesus
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Copy number formats from existing cells so new cells reuse the same
# --- style indices (date format from C2, percentage format from F2) ---
$ws.Range("C2").Copy()
$ws.Range("C7:C10").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F10").PasteSpecial(-4122)

# --- Enter the text values in the same order the cells were authored so ---
# --- the shared-strings table comes out in the expected sequence        ---
$ws.Range("A7").Value = "Repositorio listo"
$ws.Range("A9").Value = "Configuración de NetBeans para trabajar en el proyecto_1"
$ws.Range("E9").Value = "6:30pm"
$ws.Range("E7").Value = "3:00pm"
$ws.Range("A10").Value = "Configuración de NetBeans para trabajar en el proyecto_1.1"
$ws.Range("A8").Value = "Revisión_1"
$ws.Range("B8").Value = "Sirio"
$ws.Range("D10").Value = "6:00am"
$ws.Range("E10").Value = "6:30am"

# --- Remaining text values (reuse already-known shared strings) ---
$ws.Range("B7").Value = "Jesús"
$ws.Range("D7").Value = "2:00pm"
$ws.Range("B9").Value = "Jesús"
$ws.Range("D9").Value = "2:00pm"
$ws.Range("B10").Value = "Jesús"

# --- Dates ---
$ws.Range("C7").Value = (Get-Date -Year 2015 -Month 4 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C8").Value = (Get-Date -Year 2015 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C9").Value = (Get-Date -Year 2015 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C10").Value = (Get-Date -Year 2015 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0)

# --- Percentages ---
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("F9").Value = 0.9
$ws.Range("F10").Value = 1

# --- Column A is widened to fit the longer activity names ---
$ws.Columns.Item(1).ColumnWidth = 54

# --- Final selection left on F10, matching the end of the edit session ---
$ws.Range("F10").Select()
